$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 (the "נספחים" / "סקר ספרות" task): work continued into the next
#     day, so update the end date and the gross/net work-day counts.
$ws.Range("F17").Value = 43986.793055555558
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = 1.5

# --- New row 18: a new practical Kaggle task about the NYC taxi fare
#     prediction dataset.
$ws.Range("A18").Value = 4
$ws.Range("B18").Value = "מעשי"
$ws.Range("C18").Value = "Kaggle"
$ws.Range("D18").Value = "עבודה עם הדאטה של NYC taxi fare prediction"
$ws.Range("E18").Value = 43986.59097222222

# --- Fix a typo in the notes for row 16 (I16): "רלוונטים לא" -> "רלוונטים ולא רלוונטים"
$ws.Range("I16").Value = "במאמר: לפרט פרקים רלוונטים ולא רלוונטים (למשל: ניסויים פחות מעניין). `nבלוגפוסטים קשוחים (LIGHTGBM, CATBOOST).`nבלוגפוסט על CATBOOST שעוזר להבין את הרעיון של Ordered Boosting: http://towardsdatascience.com/catboost-d1f1366aca34`nיש סדרת סרטונים ביוטיוב (סה`"כ ~100 דקות) של statquest שמסבירים XGBOOST בצורה מאד פשוטה"

# --- Move the cursor/selection to B16 (matches the author's new viewport / selection position)
$ws.Range("B16").Select()
